$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (weather reading) with the new values from the diff.
$ws.Range("B2").Value = 45640.020833333336
$ws.Range("C2").Value = 45640.96666666667
$ws.Range("D2").Value = 1.6
$ws.Range("E2").Value = 7.1
$ws.Range("F2").Value = 3.73
$ws.Range("G2").Value = 2.9

# Remove the now-obsolete rows 3 and 4 (shifts dimension down to A1:G2).
$ws.Rows("3:4").Delete()
